$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46005
$ws.Range("B2").Value = 86.78
$ws.Range("C2").Value = 76.06
$ws.Range("D2").Value = 70.73999999999999
$ws.Range("E2").Value = 68.45999999999999
$ws.Range("F2").Value = 65.08
$ws.Range("G2").Value = 66.03
$ws.Range("H2").Value = 70.76000000000001
$ws.Range("I2").Value = 72.53
$ws.Range("J2").Value = 73.47
$ws.Range("K2").Value = 72.36
$ws.Range("L2").Value = 68.09999999999999
$ws.Range("M2").Value = 61.69
$ws.Range("N2").Value = 62.29
$ws.Range("O2").Value = 59
$ws.Range("P2").Value = 59.13
$ws.Range("Q2").Value = 63.42
$ws.Range("R2").Value = 74.89
$ws.Range("S2").Value = 93.13
$ws.Range("T2").Value = 105.97
$ws.Range("U2").Value = 105.94
$ws.Range("V2").Value = 106.68
$ws.Range("W2").Value = 106.06
$ws.Range("X2").Value = 102.15
$ws.Range("Y2").Value = 94.90000000000001
$ws.Range("Z2").Value = 78.56999999999999
$ws.Range("AB2").Value = 102.45
$ws.Range("AD2").Value = 106.37
$ws.Range("AF2").Value = 105.96
$ws.Range("AG2").Value = "1h-16h"
